$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2:C6 with new values
$ws.Range("B2").Value = 306
$ws.Range("C2").Value = 0.9981020682082578

$ws.Range("B3").Value = 446
$ws.Range("C3").Value = 0.9980273170769528

$ws.Range("B4").Value = 551
$ws.Range("C4").Value = 0.9980064553618823

$ws.Range("B5").Value = 803
$ws.Range("C5").Value = 0.9979077879693319

$ws.Range("B6").Value = 1005
$ws.Range("C6").Value = 0.9978000776964508

# Remove rows 7 through 16 entirely (data previously present there)
$ws.Range("A7:C16").EntireRow.Delete()
